$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly price rows (originally rows 2-4) are rotated:
#   row2 <- old row4, row3 <- old row2, row4 <- old row3.
# Capture the "before" values first (via Value2, which returns the raw
# cell value rather than a property descriptor) so the rotation is
# correct regardless of execution order.

$row2 = @{
    D = $ws.Range("D2").Value2
    J = $ws.Range("J2").Value2
    K = $ws.Range("K2").Value2
    L = $ws.Range("L2").Value2
    M = $ws.Range("M2").Value2
    P = $ws.Range("P2").Value2
}

$row3 = @{
    D = $ws.Range("D3").Value2
    J = $ws.Range("J3").Value2
    K = $ws.Range("K3").Value2
    L = $ws.Range("L3").Value2
    M = $ws.Range("M3").Value2
    P = $ws.Range("P3").Value2
}

$row4 = @{
    D = $ws.Range("D4").Value2
    J = $ws.Range("J4").Value2
    K = $ws.Range("K4").Value2
    L = $ws.Range("L4").Value2
    M = $ws.Range("M4").Value2
    P = $ws.Range("P4").Value2
}

# Row 2 takes former row 4 values
$ws.Range("D2").Value2 = $row4.D
$ws.Range("J2").Value2 = $row4.J
$ws.Range("K2").Value2 = $row4.K
$ws.Range("L2").Value2 = $row4.L
$ws.Range("M2").Value2 = $row4.M
$ws.Range("P2").Value2 = $row4.P

# Row 3 takes former row 2 values
$ws.Range("D3").Value2 = $row2.D
$ws.Range("J3").Value2 = $row2.J
$ws.Range("K3").Value2 = $row2.K
$ws.Range("L3").Value2 = $row2.L
$ws.Range("M3").Value2 = $row2.M
$ws.Range("P3").Value2 = $row2.P

# Row 4 takes former row 3 values
$ws.Range("D4").Value2 = $row3.D
$ws.Range("J4").Value2 = $row3.J
$ws.Range("K4").Value2 = $row3.K
$ws.Range("L4").Value2 = $row3.L
$ws.Range("M4").Value2 = $row3.M
$ws.Range("P4").Value2 = $row3.P

$wb.Save()
